$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: task moved out of "In Progress" (B) into "Done" (D)
$ws.Range("B2").Value = $null
$ws.Range("D2").Value = "X"

# Row 3: task moved out of "In Progress" (B) into "Done" (D)
$ws.Range("B3").Value = $null
$ws.Range("D3").Value = "X"

# Row 4: task moved out of "In Progress" (B) into "Review" (C)
$ws.Range("B4").Value = $null
$ws.Range("C4").Value = "X"

# Update the active cell selection
$ws.Range("D3").Select()
